$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "Giver the above tables solve the following queries."
#   -> "Given" (run 1, unchanged formatting)
#    + " the above tables solve the following queries." (run 2, same
#      formatting, new run boundary, leading space preserved)
# -----------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Giver the above tables solve the following queries.")
if ($found1) {
    $start = $rng1.Start

    # Fix the typo "Giver" -> "Given" in place.
    $word1 = $d.Range($start, $start + 5)
    $word1.Text = "Given"

    # Force a run split right after "Given" by toggling a character
    # format on the remainder of the sentence (collapses back to the
    # same formatting afterwards, but the run boundary now exists).
    $restLen = "the above tables solve the following queries.".Length + 1
    $rest = $d.Range($start + 5, $start + 5 + $restLen)
    $rest.Font.Bold = $false
    $rest.Font.Bold = $true
}

# -----------------------------------------------------------------
# Change 2: wrap the "ANS." run in a gramStart/gramEnd proofErr pair.
# -----------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("ANS.")
if ($found2) {
    $s = $rng2.Start
    $e = $rng2.End
    $target = $d.Range($s, $e)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype" w:cstheme="minorHAnsi"/><w:b/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>ANS.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}

Write-Output "edit complete"
